$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently has 4 metric blocks (5 rows each, rows 2-21):
#   f1_macro_mean (2-6), f1_micro_mean (7-11), f1_macro_std (12-16), f1_micro_std (17-21)
# We add a new "accuracy_balanced" metric, inserting its "mean" block right after
# the existing "mean" blocks (so before f1_macro_std) and its "std" block at the
# very end (after f1_micro_std). This matches the target layout:
#   f1_macro_mean (2-6), f1_micro_mean (7-11), accuracy_balanced_mean (12-16),
#   f1_macro_std (17-21), f1_micro_std (22-26), accuracy_balanced_std (27-31)
# ---------------------------------------------------------------------------

# Step 1: insert 5 blank rows before the current row 12 (shifts old rows 12-21 down to 17-26)
$ws.Range("A12:H16").Insert()
# Restore the bold/bordered style (column A) that Insert() does not carry over correctly;
# row 17 still has the original formatting of what used to be row 12.
$ws.Range("A17:A21").Copy()
$ws.Range("A12:A16").PasteSpecial(-4122)

# Step 2: insert 5 more blank rows after the (now shifted) f1_micro_std block, i.e. before row 27
$ws.Range("A27:H31").Insert()
$ws.Range("A22:A26").Copy()
$ws.Range("A27:A31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Fill the "n_sample" (column B) labels for both new blocks using the existing
# text values ("0","100","500","1000","2116 (all)") so the cell type stays text
# rather than being auto-converted to a number by a plain .Value assignment.
# ---------------------------------------------------------------------------
$ws.Range("B2:B6").Copy()
$ws.Range("B12").PasteSpecial(-4163)

$ws.Range("B2:B6").Copy()
$ws.Range("B27").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# accuracy_balanced_mean (rows 12-16)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "accuracy_balanced_mean"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0.459

$ws.Range("A13").Value = "accuracy_balanced_mean"
$ws.Range("C13").Value = 0.586
$ws.Range("D13").Value = 0.529
$ws.Range("E13").Value = 0.673
$ws.Range("F13").Value = 0.642
$ws.Range("G13").Value = 0.719
$ws.Range("H13").Value = 0.721

$ws.Range("A14").Value = "accuracy_balanced_mean"
$ws.Range("C14").Value = 0.635
$ws.Range("D14").Value = 0.6
$ws.Range("E14").Value = 0.735
$ws.Range("F14").Value = 0.726
$ws.Range("G14").Value = 0.795
$ws.Range("H14").Value = 0.814

$ws.Range("A15").Value = "accuracy_balanced_mean"
$ws.Range("C15").Value = 0.627
$ws.Range("D15").Value = 0.633
$ws.Range("E15").Value = 0.732
$ws.Range("F15").Value = 0.734
$ws.Range("G15").Value = 0.804
$ws.Range("H15").Value = 0.835

$ws.Range("A16").Value = "accuracy_balanced_mean"
$ws.Range("C16").Value = 0.685
$ws.Range("D16").Value = 0.688
$ws.Range("E16").Value = 0.731
$ws.Range("F16").Value = 0.755
$ws.Range("G16").Value = 0.798
$ws.Range("H16").Value = 0.841

# ---------------------------------------------------------------------------
# accuracy_balanced_std (rows 27-31)
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "accuracy_balanced_std"
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0

$ws.Range("A28").Value = "accuracy_balanced_std"
$ws.Range("C28").Value = 0.015
$ws.Range("D28").Value = 0.027
$ws.Range("E28").Value = 0.021
$ws.Range("F28").Value = 0.033
$ws.Range("G28").Value = 0.027
$ws.Range("H28").Value = 0.034

$ws.Range("A29").Value = "accuracy_balanced_std"
$ws.Range("C29").Value = 0.015
$ws.Range("D29").Value = 0.007
$ws.Range("E29").Value = 0.012
$ws.Range("F29").Value = 0.015
$ws.Range("G29").Value = 0.005
$ws.Range("H29").Value = 0.017

$ws.Range("A30").Value = "accuracy_balanced_std"
$ws.Range("C30").Value = 0.001
$ws.Range("D30").Value = 0.008
$ws.Range("E30").Value = 0.02
$ws.Range("F30").Value = 0.011
$ws.Range("G30").Value = 0.014
$ws.Range("H30").Value = 0.009

$ws.Range("A31").Value = "accuracy_balanced_std"
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0.002
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0.011
$ws.Range("H31").Value = 0.006

Write-Host ("Final dimension: {0}" -f $ws.UsedRange.Address())
